$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 71, pushing the existing rows 71-163 down to 72-164.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new data record.
$ws.Range("A71").Value = 5
$ws.Range("B71").Value = "Macroferia Regional de Talca"
$ws.Range("C71").Value = "Maule"
$ws.Range("D71").Value = 44413
$ws.Range("E71").Value = 7
$ws.Range("F71").Value = 100112032
$ws.Range("G71").Value = "Zapallo italiano"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 300
$ws.Range("K71").Value = 7000
$ws.Range("L71").Value = 7000
$ws.Range("M71").Value = 7000
$ws.Range("N71").Value = "`$/caja 50 unidades"
$ws.Range("O71").Value = "Región del Maule"
$ws.Range("P71").Value = 140
$ws.Range("Q71").Value = 50
$ws.Range("R71").Value = "Hortaliza"
